$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    # Round-trip through a scratch formula cell + copy/paste-values so the
    # destination ends up as a genuine shared-string text cell with the
    # *default* style (no stray NumberFormat / quote-prefix residue).
    $ws.Range("Z100").Formula = '="' + $text + '"'
    $ws.Range("Z100").Copy()
    $ws.Range($addr).PasteSpecial(-4163)
    $ws.Range("Z100").Clear()
}

function Set-EmptyText($addr) {
    # Force a real empty-string text cell (t="s" -> "") rather than letting
    # the engine treat an assigned "" as "clear the cell".
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = "'"
    $ws.Range($addr).ClearFormats()
}

# ---------------------------------------------------------------------
# Row 2 edits
# ---------------------------------------------------------------------
Set-TextValue "A2" "06250006"
$ws.Range("B2").Value = 45820.34789040509
$ws.Range("C2").Value = "--"
Set-TextValue "D2" "06250001"
$ws.Range("E2").Value = "--"
$ws.Range("F2").Value = "SALLE D'ATTENTE NIVEAU 1"
$ws.Range("G2").Clear()
$ws.Range("H2").Clear()
$ws.Range("K2").Clear()
$ws.Range("L2").Value = "N/C"

# ---------------------------------------------------------------------
# Row 3 edits
# ---------------------------------------------------------------------
Set-TextValue "A3" "06250007"
$ws.Range("B3").Value = 45820.63681650463
$ws.Range("C3").Value = "--"
$ws.Range("D3").Value = "--"
$ws.Range("E3").Value = "--"
$ws.Range("F3").Value = "SALLE DE CONFERENCE ROOM"
Set-EmptyText "I3"
$ws.Range("J3").Clear()
$ws.Range("K3").Clear()
$ws.Range("L3").Value = 45820.64248148148

# ---------------------------------------------------------------------
# Row 4 (new row)
# ---------------------------------------------------------------------
Set-TextValue "A4" "06250003"
$ws.Range("B4").Value = 45819.60005408565
$ws.Range("B2").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("C4").Value = "--"
$ws.Range("D4").Value = "--"
$ws.Range("E4").Value = "--"
$ws.Range("F4").Value = "SALLE DE CONFERENCE ROOM"
Set-EmptyText "I4"
$ws.Range("L4").Value = "N/C"
$ws.Range("M4").Value = "Admin User"
$ws.Range("O4").Value = "EN ATTENTE"
